$ErrorActionPreference = "Stop"

function Get-ParaByText($pres, $needle) {
    for ($si = 1; $si -le $pres.Slides.Count; $si++) {
        $slide = $pres.Slides.Item($si)
        for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
            $shape = $slide.Shapes.Item($shi)
            if ($shape.HasTextFrame -eq $false) { continue }
            $tf = $shape.TextFrame
            if ($tf.HasText -eq $false) { continue }
            $tr = $tf.TextRange
            $count = $tr.Paragraphs().Count
            for ($pi = 1; $pi -le $count; $pi++) {
                $para = $tr.Paragraphs($pi)
                if ($para.Text -like "*$needle*") {
                    return $para
                }
            }
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Edit 1 (slide "First Application"): merge the three runs of the first
# content paragraph into a single run reading:
# "Now that we have all our development tools in place, let's create our
# first android application in Android Studio."
# ---------------------------------------------------------------------------
$para1 = Get-ParaByText $p "we have all our development tools"

$mergedText = "Now that we have all our development tools in place, let" + [char]0x2019 + "s create our first android application in Android Studio."
$fullRange1 = $para1.Characters(1, $para1.Text.Length)
$fullRange1.Text = $mergedText

# ---------------------------------------------------------------------------
# Edit 2 (slide "Assignment"): split the single run
# "(In not less than 200 words)" into three runs and change 200 -> 100:
#   "(In not less "  /  "than 100 "  /  "words)"
# ---------------------------------------------------------------------------
$para2 = Get-ParaByText $p "In not less than"

$part1 = "(In not less "
$part2 = "than 100 "
$part3 = "words)"

$r1 = $para2.Characters(1, $part1.Length)
$r1.Text = $part1

$r2 = $para2.Characters($part1.Length + 1, $part2.Length)
$r2.Text = $part2

$r3 = $para2.Characters($part1.Length + $part2.Length + 1, $part3.Length)
$r3.Text = $part3
